$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.687.71'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '2.173.75'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '238.35'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("D7").Value = '72.48'
$ws.Range("E7").Value = '  -2.99%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("D10").Value = '40.50'
$ws.Range("E10").Value = '  -4.56%  '
$ws.Range("D11").Value = '0.0911'
$ws.Range("E11").Value = '  -5.07%  '
$ws.Range("D12").Value = '54.69'
$ws.Range("E12").Value = '  -3.25%  '
$ws.Range("D13").Value = '6.75'
$ws.Range("E13").Value = '  -3.06%  '
$ws.Range("D14").Value = '0.0999'
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("D15").Value = '2.500.24'
$ws.Range("E15").Value = '  -2.75%  '
$ws.Range("D16").Value = '14.42'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '2.175.50'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("E18").Value = '  -6.88%  '
$ws.Range("D19").Value = '41.635.29'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '70.19'
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("D22").Value = '5.81'
$ws.Range("E22").Value = '  -6.79%  '
$ws.Range("D23").Value = '10.02'
$ws.Range("E23").Value = '  -12.03%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").Value = '10.76'
$ws.Range("E27").Value = '  -5.30%  '
$ws.Range("E28").Value = '  -9.96%  '
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  -3.67%  '
$ws.Range("D30").Value = '170.92'
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("E31").Value = '  -5.48%  '
$ws.Range("E32").Value = '  -3.67%  '
$ws.Range("D33").Value = '32.94'
$ws.Range("E33").Value = '  +10.20%  '
$ws.Range("D34").Value = '0.0778'
$ws.Range("E34").Value = '  -3.37%  '
$ws.Range("D35").Value = '5.30'
$ws.Range("D36").Value = '0.121'
$ws.Range("E36").Value = '  -3.44%  '
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  -5.80%  '
$ws.Range("D39").Value = '0.0313'
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '12.13'
$ws.Range("E40").Value = '  -8.39%  '
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("E42").Value = '  -5.89%  '
$ws.Range("D43").Value = '59.23'
$ws.Range("E43").Value = '  -8.89%  '
$ws.Range("D44").Value = '0.191'
$ws.Range("E44").Value = '  -4.66%  '
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").Value = '97.64'
$ws.Range("E47").Value = '  -6.75%  '
$ws.Range("E48").Value = '  -4.24%  '
$ws.Range("E49").Value = '  -4.83%  '
$ws.Range("E50").Value = '  -6.43%  '
$ws.Range("E51").Value = '  -2.14%  '
